$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Eggers, Vivyan, and Wagner" / "Eggers et al." study row (row 15),
# shifting all subsequent rows up by one (full replication of this study's
# results was incorporated, replacing the placeholder entry).
$ws.Rows.Item(15).Delete()

# Leave the active cell/selection where the editor last left it.
[void]$ws.Range("C12").Select()
